$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 16

$ws.Cells.Item($row, 1).Value = 112552453          # A - Id
$ws.Cells.Item($row, 2).Value = 77650               # B - Taxonsorteringsordning
$ws.Cells.Item($row, 3).Value = "Ovaliderad"        # C - Valideringsstatus
$ws.Cells.Item($row, 4).Value = "NT"                # D - Rödlistade
$ws.Cells.Item($row, 5).Value = 6425                # E - TaxonId
$ws.Cells.Item($row, 6).Value = "Garnlav"           # F - Artnamn
$ws.Cells.Item($row, 7).Value = "Alectoria sarmentosa"  # G - Vetenskapligt namn
$ws.Cells.Item($row, 8).Value = "(Ach.) Ach."        # H - Auktor

$ws.Cells.Item($row, 16).Value = "Sydväst Skarptjärnen, Vrm"  # P - Lokalnamn
$ws.Cells.Item($row, 17).Value = 404648              # Q - Ost
$ws.Cells.Item($row, 18).Value = 6706779             # R - Nord
$ws.Cells.Item($row, 19).Value = 10                  # S - Noggrannhet
$ws.Cells.Item($row, 20).Value = "Värmland"          # T - Län
$ws.Cells.Item($row, 21).Value = "Torsby"            # U - Kommun
$ws.Cells.Item($row, 22).Value = "Värmland"          # V - Provins
$ws.Cells.Item($row, 23).Value = "Norra Ny"          # W - Församling

$ws.Cells.Item($row, 25).Value = "'2023-07-27"       # Y - Startdatum (force text, avoid date parsing)
$ws.Cells.Item($row, 27).Value = "'2023-07-27"       # AA - Slutdatum (force text, avoid date parsing)

$ws.Cells.Item($row, 30).Value = $false              # AD - Ej återfunnen
$ws.Cells.Item($row, 31).Value = $false              # AE - Osäker artbestämning
$ws.Cells.Item($row, 33).Value = $false              # AG - Ospontan

$ws.Cells.Item($row, 49).Value = "Anders Boström"    # AW - Rapportör
$ws.Cells.Item($row, 50).Value = "Anders Boström"    # AX - Observatörer
